# Edit script: apply documentation updates to Word doc via COM automation
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the "Fig5. Database Design" image run as NoProofing (adds <w:noProof/>)
#    This is paragraph 46 (1-indexed) -- the image-only paragraph right after
#    the "Database design" heading.
# ---------------------------------------------------------------------------
$imgPara = $d.Paragraphs.Item(46)
$imgPara.Range.NoProofing = $true

# ---------------------------------------------------------------------------
# 2) Insert a new bullet under the Question("/questions") section describing
#    the "/delete" endpoint, right after the existing "/vote" bullet (item 64)
#    and before the "Answer("/answers")" heading (item 65).
# ---------------------------------------------------------------------------
$questionVotePara = $d.Paragraphs.Item(64)
$insertAfterQuestionVote = $d.Range($questionVotePara.Range.End - 1, $questionVotePara.Range.End - 1)
$questionDeleteXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listparagraf"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/delete</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">” – Query Params is </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>questionid</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, the id of the question that is to be deleted.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertAfterQuestionVote.InsertXML($questionDeleteXml)

# ---------------------------------------------------------------------------
# 3) Move the rendered-page-break marker: it currently sits on the "/update"
#    bullet of the Answer("/answers") section (now item 69, after the insert
#    above shifted everything below it down by one); it should instead sit on
#    the preceding "/save" bullet (item 68).
# ---------------------------------------------------------------------------
$answerSavePara = $d.Paragraphs.Item(68)
$answerUpdatePara = $d.Paragraphs.Item(69)

# 3a. Remove it from "/update"'s opening quote run.
$updateStart = $answerUpdatePara.Range.Start
$updateQuoteRange = $d.Range($updateStart, $updateStart + 1)
$updateQuoteRange.Delete()
$updateInsertPoint = $d.Range($updateStart, $updateStart)
$updateQuoteXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>“</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$updateInsertPoint.InsertXML($updateQuoteXml)

# 3b. Add it to "/save"'s opening quote run.
$saveStart = $answerSavePara.Range.Start
$saveQuoteRange = $d.Range($saveStart, $saveStart + 1)
$saveQuoteRange.Delete()
$saveInsertPoint = $d.Range($saveStart, $saveStart)
$saveQuoteXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>“</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$saveInsertPoint.InsertXML($saveQuoteXml)

# ---------------------------------------------------------------------------
# 4) Append the new Tag("/tags") endpoint section at the very end of the
#    document (after the last "/vote" bullet of Answer("/answers"), before
#    the section properties).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfDoc = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$tagSectionXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listparagraf"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Tag(“</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/tags</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>”)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listparagraf"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/getAll</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">” – No body, returns all </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tags</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> in the database.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listparagraf"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/save</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">” – Body is a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tag</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, containing </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listparagraf"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/update</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">” – Query Params </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tagid</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">the body contains the new </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> of tag.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listparagraf"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/delete</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">” – Query Params has the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tagid</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, the id of the </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tag</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to be deleted.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$endOfDoc.InsertXML($tagSectionXml)

Write-Output "edit complete"
